$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "sample" value for row 7 (Sample Type column), matching the
# existing "sample"/"blank" entries used elsewhere in the sheet
$ws.Range("B7").Value = "sample"

# Move the active selection to B7 (also resets the view's scrolled
# topLeftCell back to the sheet's natural origin)
$ws.Range("B7").Select()

$wb.Save()
